$d = $word.ActiveDocument

# 1) Remove the trailing period after "order placement time from 9 to 2 seconds."
$d.Content.Find.Execute(
    "order placement time from 9 to 2 seconds.", $true, $false, $false, $false, $false,
    $true, 1, $false, "order placement time from 9 to 2 seconds", 2) | Out-Null

# 2) Remove the trailing period after "...the application gained 30% in throughput."
$d.Content.Find.Execute(
    "the application gained 30% in throughput.", $true, $false, $false, $false, $false,
    $true, 1, $false, "the application gained 30% in throughput", 2) | Out-Null

# 3) Remove the trailing period after "...3,000 hits per second generated traffic."
$d.Content.Find.Execute(
    "3,000 hits per second generated traffic.", $true, $false, $false, $false, $false,
    $true, 1, $false, "3,000 hits per second generated traffic", 2) | Out-Null

# 4) Move the "_GoBack" bookmark from the header row (after the 3rd tab run) to the
#    end of the paragraph that now ends in "...generated traffic" (no period), in the
#    last table cell. Adding a bookmark named "_GoBack" again removes the old one
#    (bookmark names are unique) and creates the new one at the given range.
$r = $d.Content
$r.Find.Execute("ava clients for high load emulation, so it allowed to verify back-end services under 3,000 hits per second generated traffic") | Out-Null
$bmRange = $r.Paragraphs(1).Range
$bmRange.Collapse(0)
$bmRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
